{"js": "const body = context.document.body;\nconst results = body.search(\"1/30/2024 9:11 AM\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\nif (results.items.length > 0) {\n  results.items[0].insertText(\"2/1/2024 1:12 PM\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$rng = $d.Content\n$rng.Find.Execute(\"1/30/2024 9:11 AM\", $false, $false, $false, $false, $false, $true, 1, $false, \"2/1/2024 1:12 PM\", 2)\n"}
